$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 668
$ws.Range("J3").Value = 734
$ws.Range("I4").Value = 1753
$ws.Range("J4").Value = 160
$ws.Range("J5").Value = 52
$ws.Range("J6").Value = 1069
$ws.Range("I7").Value = 26184
$ws.Range("J7").Value = 2683

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 8
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 101

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 13
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 19
$ws.Range("I8").Value = 1543
$ws.Range("J8").Value = 169
$ws.Range("J15").Value = 32
$ws.Range("J20").Value = 50
$ws.Range("J25").Value = 16
$ws.Range("J27").Value = 11
$ws.Range("J29").Value = 136
$ws.Range("J33").Value = 111
$ws.Range("J36").Value = 41
$ws.Range("J37").Value = 101
$ws.Range("J42").Value = 116
$ws.Range("J43").Value = 37
$ws.Range("J44").Value = 22
$ws.Range("J49").Value = 13
$ws.Range("J51").Value = 34
$ws.Range("J53").Value = 28
$ws.Range("J54").Value = 45
$ws.Range("J60").Value = 16
$ws.Range("J63").Value = 11
$ws.Range("J64").Value = 17
$ws.Range("J66").Value = 6
$ws.Range("J67").Value = 100
$ws.Range("J68").Value = 7
$ws.Range("J73").Value = 27
$ws.Range("J75").Value = 11
$ws.Range("J78").Value = 31
$ws.Range("J79").Value = 89
$ws.Range("J82").Value = 6
$ws.Range("J85").Value = 110
$ws.Range("J88").Value = 21
$ws.Range("J89").Value = 30
$ws.Range("J91").Value = 33
$ws.Range("J94").Value = 16
$ws.Range("J95").Value = 51
$ws.Range("J96").Value = 33
$ws.Range("J97").Value = 16
$ws.Range("J99").Value = 36
$ws.Range("I101").Value = 26184
$ws.Range("J101").Value = 2683

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 26
$ws.Range("J3").Value = 28
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 13
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 43
$ws.Range("J3").Value = 47
$ws.Range("J7").Value = 136

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 5
$ws.Range("J3").Value = 10

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 26
$ws.Range("J3").Value = 35
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 16
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 26
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 7
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 9
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J4").Value = 1
$ws.Range("J7").Value = 6

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 54
$ws.Range("J3").Value = 56
$ws.Range("I5").Value = 48
$ws.Range("I7").Value = 1543
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 7

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 3
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("J5").Value = 2
$ws.Range("J6").Value = 6

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 20
$ws.Range("J3").Value = 33
